$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Copy the formatting of the last existing data row (row 9) down onto
#    the three new rows (10-12) before we touch row 9's own styling, so
#    the engine's style table reuses the same "shape" of xf entries
#    (including the built-in date number format already used by column D).
# ---------------------------------------------------------------------
$ws.Range("A9:J9").Copy()
$ws.Range("A10:J12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Extend the AutoFilter range from A1:H8 to A1:H9 *before* filling in
#    the new rows, so the engine doesn't silently grow the filter range
#    to match the newly-used range. Keep the workbook-level
#    _FilterDatabase defined name in sync too.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:H9").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$9"
    }
}

# ---------------------------------------------------------------------
# 3) Row 10 - LeetCode 94, Binary Tree Inorder Traversal
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "LeetCode"
$ws.Range("B10").Value = 94
$ws.Range("C10").Value = "二叉树的中序遍历"
$ws.Range("D10").Value = "8/24/2022"
$ws.Range("E10").Value = "0星"
$ws.Range("F10").Value = "dfs,二叉树中序遍历"
$ws.Range("G10").Value = "给出二叉树，输出其中序遍历"
$ws.Range("H10").Value = "dfs先左子树，再跟节点，后右子树，记录中序遍历最终打印"
$ws.Range("I10").Value = "Wyatt"

# C10 also gets word-wrap turned on (new style compared to the other cells)
$ws.Range("C10").WrapText = $true

# ---------------------------------------------------------------------
# 4) Row 11 - LeetCode 58, Length of Last Word
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "LeetCode"
$ws.Range("B11").Value = 58
$ws.Range("C11").Value = "最后一个单词的长度"
$ws.Range("D11").Value = "8/24/2022"
$ws.Range("E11").Value = "0星"
$ws.Range("F11").Value = "——"
$ws.Range("G11").Value = "给出只包含单词和空格的字符串，求最后一个单词的长度"
$ws.Range("H11").Value = "简单字符串处理"
# I11 (submitter) intentionally left blank, but still bordered/centered
$ws.Range("I11").Borders.LineStyle = 1
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("I11").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5) Row 12 - LeetCode 1460, Make Two Arrays Equal by Reversing Sub-arrays
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "LeetCode"
$ws.Range("B12").Value = 1460
$ws.Range("C12").Value = "通过翻转子数组使两个数组相等"
$ws.Range("D12").Value = "8/24/2022"
$ws.Range("E12").Value = "0星"
$ws.Range("F12").Value = "——"
$ws.Range("G12").Value = "给出两个等长数组，判断是否能通过若干次子数组反转使得二者一样"
$ws.Range("H12").Value = "冒泡排序就是多次的长度为2的子数组反转，所以只要两个数组内的相应元素个数相等即一定可以达成"
# I12 (submitter) intentionally left blank, but still bordered/centered
$ws.Range("I12").Borders.LineStyle = 1
$ws.Range("I12").HorizontalAlignment = -4108
$ws.Range("I12").VerticalAlignment = -4108
$ws.Range("C12").WrapText = $true

# ---------------------------------------------------------------------
# 6) Normalise row 9's own styling now that the new rows exist (collapses
#    the duplicate xf entries this sheet had accumulated back to the
#    canonical ones used by rows 2-8).
# ---------------------------------------------------------------------
$ws.Range("A9:J9").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 7) Extend the star-rating data validation list (column E) to cover the
#    new rows too.
# ---------------------------------------------------------------------
$ws.Range("E2:E12").Validation.Delete()
$ws.Range("E2:E12").Validation.Add(3, 1, 1, '"0星,1星,2星,3星,4星,5星"')

# ---------------------------------------------------------------------
# 8) Match the final selection left behind in the source workbook.
# ---------------------------------------------------------------------
$ws.Range("H12").Select()
